$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Free_ligand_(uM)"
$ws.Range("B1").Value = "n_bar"

$ws.Range("E6").Select()
